$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.829.15'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.502.25'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -4.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.16'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.03'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.520'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.501.90'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.62%  '
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("E12").Value = '  -4.25%  '
$ws.Range("E13").Value = '  -2.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.19'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -5.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.940.40'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.26%  '
$ws.Range("E16").Value = '  -4.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.593.93'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.508.41'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -4.36%  '
$ws.Range("E19").Value = '  -7.35%  '
$ws.Range("E20").Value = '  -4.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.30'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.95%  '
$ws.Range("E22").Value = '  -3.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.53'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.00%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.92'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").Value = '  -1.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.86'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.626.84'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0972'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.16%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '518.91'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.70%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.02'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.09%  '
$ws.Range("E33").Value = '  -3.33%  '
$ws.Range("E34").Value = '  -5.28%  '
$ws.Range("E35").Value = '  -4.17%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '156.82'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("E38").Value = '  -4.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.44'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.98%  '
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.353'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.69%  '
$ws.Range("E42").Value = '  -3.41%  '
$ws.Range("E43").Value = '  -3.84%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '146.09'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.00%  '
$ws.Range("E47").Value = '  -4.80%  '
$ws.Range("E48").Value = '  -2.83%  '
$ws.Range("E49").Value = '  -8.35%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0748'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.86%  '
